# Rename "portfolio" to "stocks", add a new "portfolio" summary sheet after it,
# and populate the new sheet with an asset/value table (CASH + STOCKS rows).

$wb = $excel.ActiveWorkbook

# 1. Rename the existing "portfolio" sheet (2nd sheet) to "stocks"
$stocks = $wb.Worksheets.Item(2)
$stocks.Name = "stocks"

# 2. Widen column B on the stocks sheet (ticker values need more room)
$stocks.Columns.Item(2).ColumnWidth = 14.86

# 3. Add a brand-new "portfolio" sheet right after "stocks"
$portfolio = $wb.Worksheets.Add($null, $stocks)
$portfolio.Name = "portfolio"

# 4. Populate the new portfolio sheet (order matters for shared-string ids)
$portfolio.Range("A2").Value = "CASH"
$portfolio.Range("A1").Value = "asset"
$portfolio.Range("B1").Value = "value"
$portfolio.Range("B2").Value = 10000
$portfolio.Range("A3").Value = "STOCKS"

# 5. Update selections to match the final state
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("E2").Select()

[void]$portfolio.Range("C1").Select()
